# Fig_S2.xlsx — "Add files via upload" edit
#
# The source workbook had 9 sheets named after 10-year age brackets
# ("0-9", "10-19", ..., "80 or above"), each holding the same F-column
# layout of data used to build Figure S2 (panels a-i). The re-upload
# renames every sheet to the matching figure-panel label (Fig_S2a..i)
# and leaves the active/selected sheet on the last tab (Fig_S2i) instead
# of the first one.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "Fig_S2a",
    "Fig_S2b",
    "Fig_S2c",
    "Fig_S2d",
    "Fig_S2e",
    "Fig_S2f",
    "Fig_S2g",
    "Fig_S2h",
    "Fig_S2i"
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $sheet = $wb.Worksheets.Item($i + 1)
    $sheet.Name = $newNames[$i]
}

# Make the last sheet (Fig_S2i) the active/selected tab, matching the
# workbook's saved view state (tabSelected moves off the first sheet).
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
